$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145, pushing the existing rows 145:211 down to 146:212.
$ws.Rows("145:145").Insert()

# Populate the newly inserted row 145 with the new weekly price entry.
$ws.Range("A145").Value = 10
$ws.Range("B145").Value = "Vega Modelo de Temuco"
$ws.Range("C145").Value = "La Araucanía"
$ws.Range("D145").Value = 44466
$ws.Range("E145").Value = 9
$ws.Range("F145").Value = "Fruta"
$ws.Range("G145").Value = 100108
$ws.Range("H145").Value = "Tropicales y subtropicales"
$ws.Range("I145").Value = 100108002
$ws.Range("J145").Value = "Mango"
$ws.Range("K145").Value = "Sin especificar"
$ws.Range("L145").Value = "Primera"
$ws.Range("M145").Value = 500
$ws.Range("N145").Value = 8000
$ws.Range("O145").Value = 8000
$ws.Range("P145").Value = 8000
$ws.Range("Q145").Value = "$/bandeja 4 kilos"
$ws.Range("R145").Value = "Brasil"
$ws.Range("S145").Value = 2000
$ws.Range("T145").Value = 4
